# 25/05- Added outlier check and the Price model works as intended
# Adds a "2021" column (P) to the Table8 table, duplicating the 2020 (O)
# column's values/format as an outlier-check column, and nudges the
# active selection the way the author's session left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grow the table first so the engine has a real "column 16" to target;
# the header text gets synced into the table's column name afterwards.
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:P10"))

# Header (P1): copy O1's look (bold-ish header font/border formatting)
# then give it the literal text "2021" — the leading apostrophe forces
# it to land as a shared string instead of being auto-coerced to a
# number, matching the header cells used for the other year columns.
$ws.Range("O1").Copy($ws.Range("P1"))
$ws.Range("P1").Value = "'2021"

# Data rows (P2:P10): this is a straight copy of the 2020 (O) column,
# values and number formatting both, used as the new outlier-check
# column.
$ws.Range("O2:O10").Copy($ws.Range("P2:P10"))

# Tidy up the now-narrower "2020"/"2021" columns and give the new one a
# sensible width.
$ws.Columns.Item(14).ColumnWidth = 14.3854166666667
$ws.Columns.Item(15).ColumnWidth = 13.7213541666667
$ws.Columns.Item(16).ColumnWidth = 17.1666666666667

# Match the author's final selection in the saved session.
$ws.Range("A10").Select()
